# Scheduled market-data refresh for the Gungnir Profits workbook.
# Updates the currentAveragePrice / LevePrice / LeveProfit columns (H:N)
# on each crafting-job sheet with freshly pulled Universalis price data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 236.11111
$ws.Range("I6").Value = 132.14285
$ws.Range("K6").Value = 396.42855
$ws.Range("M6").Value = -284.42855

$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").ClearContents()

$ws.Range("H32").Value = 25000818
$ws.Range("I32").Value = 62500300
$ws.Range("J32").Value = 1163.3334
$ws.Range("K32").Value = 62500300
$ws.Range("L32").Value = 1163.3334
$ws.Range("M32").Value = -62499974
$ws.Range("N32").Value = -1815.3334

$ws.Range("H132").Value = 6332761
$ws.Range("I132").Value = 7045423.5
$ws.Range("J132").Value = 7882.375
$ws.Range("K132").Value = 21136270.5
$ws.Range("L132").Value = 23647.125
$ws.Range("M132").Value = -21133740.5
$ws.Range("N132").Value = -28707.125

$ws.Range("H138").Value = 1597.1621
$ws.Range("I138").Value = 775.625
$ws.Range("J138").Value = 3113.8462
$ws.Range("K138").Value = 2326.875
$ws.Range("L138").Value = 9341.5386
$ws.Range("M138").Value = 2813.125
$ws.Range("N138").Value = -19621.5386

$ws.Range("H141").Value = 3347.7354
$ws.Range("I141").Value = 1748.8
$ws.Range("J141").Value = 7789.222
$ws.Range("K141").Value = 5246.4
$ws.Range("L141").Value = 23367.666
$ws.Range("M141").Value = -66.39999999999964
$ws.Range("N141").Value = -33727.666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1163.9231
$ws.Range("I74").Value = 1028.7826
$ws.Range("J74").Value = 2200
$ws.Range("K74").Value = 1028.7826
$ws.Range("L74").Value = 2200
$ws.Range("M74").Value = -154.7826
$ws.Range("N74").Value = -3948

$ws.Range("H77").Value = 1163.9231
$ws.Range("I77").Value = 1028.7826
$ws.Range("J77").Value = 2200
$ws.Range("K77").Value = 5143.913
$ws.Range("L77").Value = 11000
$ws.Range("M77").Value = -775.9130000000005
$ws.Range("N77").Value = -19736

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2745.2666
$ws.Range("I20").Value = 2323.7273
$ws.Range("J20").Value = 3904.5
$ws.Range("K20").Value = 2323.7273
$ws.Range("L20").Value = 3904.5
$ws.Range("M20").Value = -2076.7273
$ws.Range("N20").Value = -4398.5

$ws.Range("H107").Value = 100001200
$ws.Range("I107").Value = 166667330
$ws.Range("J107").Value = 1996.5
$ws.Range("K107").Value = 166667330
$ws.Range("L107").Value = 1996.5
$ws.Range("M107").Value = -166665410
$ws.Range("N107").Value = -5836.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 33626.832
$ws.Range("I22").Value = 25190.25
$ws.Range("J22").Value = 50500
$ws.Range("K22").Value = 25190.25
$ws.Range("L22").Value = 50500
$ws.Range("M22").Value = -24840.25
$ws.Range("N22").Value = -51200

$ws.Range("H58").Value = 22727892
$ws.Range("I58").Value = 40000536
$ws.Range("J58").Value = 728.5263
$ws.Range("K58").Value = 40000536
$ws.Range("L58").Value = 728.5263
$ws.Range("M58").Value = -40000333
$ws.Range("N58").Value = -1134.5263

$ws.Range("H132").Value = 11112469
$ws.Range("I132").Value = 1056.2273
$ws.Range("J132").Value = 41668856
$ws.Range("K132").Value = 3168.6819
$ws.Range("L132").Value = 125006568
$ws.Range("M132").Value = -638.6819
$ws.Range("N132").Value = -125011628

$ws.Range("H134").Value = 745.7778
$ws.Range("I134").Value = 722.36365
$ws.Range("J134").Value = 848.8
$ws.Range("K134").Value = 2167.09095
$ws.Range("L134").Value = 2546.4
$ws.Range("M134").Value = 367.9090500000002
$ws.Range("N134").Value = -7616.4

$ws.Range("H136").Value = 22727892
$ws.Range("I136").Value = 40000536
$ws.Range("J136").Value = 728.5263
$ws.Range("K136").Value = 120001608
$ws.Range("L136").Value = 2185.5789
$ws.Range("M136").Value = -119999058
$ws.Range("N136").Value = -7285.5789

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 691.6667
$ws.Range("I34").Value = 383.33334
$ws.Range("J34").Value = 1000
$ws.Range("K34").Value = 1150.00002
$ws.Range("L34").Value = 3000
$ws.Range("M34").Value = -1066.00002
$ws.Range("N34").Value = -3168

$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()

$ws.Range("H55").Value = 2000
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 2000
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 6000
$ws.Range("N55").Value = -6354
$ws.Range("M55").ClearContents()

$ws.Range("H130").Value = 2147.1428
$ws.Range("J130").Value = 2416.6667
$ws.Range("L130").Value = 7250.000100000001
$ws.Range("N130").Value = -17290.0001

$ws.Range("H131").Value = 895.37
$ws.Range("I131").Value = 532.8570999999999
$ws.Range("J131").Value = 922.6559
$ws.Range("K131").Value = 1598.5713
$ws.Range("L131").Value = 2767.9677
$ws.Range("M131").Value = 3441.4287
$ws.Range("N131").Value = -12847.9677

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("M9").ClearContents()

$ws.Range("H70").Value = 4650.654
$ws.Range("I70").Value = 4541.3335
$ws.Range("J70").Value = 4799.727
$ws.Range("K70").Value = 4541.3335
$ws.Range("L70").Value = 4799.727
$ws.Range("M70").Value = -4271.3335
$ws.Range("N70").Value = -5339.727

$ws.Range("H73").Value = 4650.654
$ws.Range("I73").Value = 4541.3335
$ws.Range("J73").Value = 4799.727
$ws.Range("K73").Value = 4541.3335
$ws.Range("L73").Value = 4799.727
$ws.Range("M73").Value = -3605.3335
$ws.Range("N73").Value = -6671.727

$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("M97").ClearContents()
$ws.Range("N97").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1447.5294
$ws.Range("I68").Value = 1412
$ws.Range("J68").Value = 1487.5
$ws.Range("K68").Value = 1412
$ws.Range("L68").Value = 1487.5
$ws.Range("M68").Value = -663
$ws.Range("N68").Value = -2985.5

$ws.Range("H71").Value = 1447.5294
$ws.Range("I71").Value = 1412
$ws.Range("J71").Value = 1487.5
$ws.Range("K71").Value = 7060
$ws.Range("L71").Value = 7437.5
$ws.Range("M71").Value = -3316
$ws.Range("N71").Value = -14925.5

$ws.Range("H82").Value = 1400
$ws.Range("I82").Value = 1359.8
$ws.Range("J82").Value = 1500.5
$ws.Range("K82").Value = 1359.8
$ws.Range("L82").Value = 1500.5
$ws.Range("M82").Value = -998.8
$ws.Range("N82").Value = -2222.5

$ws.Range("H85").Value = 1400
$ws.Range("I85").Value = 1359.8
$ws.Range("J85").Value = 1500.5
$ws.Range("K85").Value = 1359.8
$ws.Range("L85").Value = 1500.5
$ws.Range("M85").Value = -111.8
$ws.Range("N85").Value = -3996.5

$ws.Range("H132").Value = 7137.475
$ws.Range("I132").Value = 2073.5264
$ws.Range("J132").Value = 11719.143
$ws.Range("K132").Value = 6220.5792
$ws.Range("L132").Value = 35157.429
$ws.Range("M132").Value = -3690.5792
$ws.Range("N132").Value = -40217.429

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H12").Value = 40573.5
$ws.Range("I12").Value = 1140
$ws.Range("K12").Value = 1140
$ws.Range("M12").Value = -998

$ws.Range("H22").Value = 9000
$ws.Range("J22").Value = 9000
$ws.Range("L22").Value = 9000
$ws.Range("N22").Value = -9586

$ws.Range("H132").Value = 32044.764
$ws.Range("I132").Value = 41800.5
$ws.Range("J132").Value = 10907.333
$ws.Range("K132").Value = 125401.5
$ws.Range("L132").Value = 32721.999
$ws.Range("M132").Value = -122871.5
$ws.Range("N132").Value = -37781.999

$ws.Range("H136").Value = 1884.2333
$ws.Range("I136").Value = 1186.9474
$ws.Range("J136").Value = 3088.6365
$ws.Range("K136").Value = 3560.8422
$ws.Range("L136").Value = 9265.9095
$ws.Range("M136").Value = -1010.8422
$ws.Range("N136").Value = -14365.9095
